$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 576/577 (pushes existing rows 576-644 down to 578-646)
$ws.Range("A576:A577").EntireRow.Insert()

# --- Row 576 ---
$ws.Cells.Item(576, 1).Value = 6
$ws.Cells.Item(576, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(576, 3).Value = "Metropolitana"
$ws.Cells.Item(576, 4).Value = 45127
$ws.Cells.Item(576, 5).Value = 13
$ws.Cells.Item(576, 6).Value = 100112043
$ws.Cells.Item(576, 7).Value = "Pepino ensalada"
$ws.Cells.Item(576, 8).Value = "Sin especificar"
$ws.Cells.Item(576, 9).Value = "Primera"
$ws.Cells.Item(576, 10).Value = 380
$ws.Cells.Item(576, 11).Value = 10000
$ws.Cells.Item(576, 12).Value = 12000
$ws.Cells.Item(576, 13).Value = 11211
$ws.Cells.Item(576, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(576, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(576, 16).Value = 187
$ws.Cells.Item(576, 17).Value = 60
$ws.Cells.Item(576, 18).Value = "Hortaliza"

# --- Row 577 ---
$ws.Cells.Item(577, 1).Value = 6
$ws.Cells.Item(577, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(577, 3).Value = "Metropolitana"
$ws.Cells.Item(577, 4).Value = 45127
$ws.Cells.Item(577, 5).Value = 13
$ws.Cells.Item(577, 6).Value = 100112043
$ws.Cells.Item(577, 7).Value = "Pepino ensalada"
$ws.Cells.Item(577, 8).Value = "Sin especificar"
$ws.Cells.Item(577, 9).Value = "Segunda"
$ws.Cells.Item(577, 10).Value = 200
$ws.Cells.Item(577, 11).Value = 10000
$ws.Cells.Item(577, 12).Value = 11000
$ws.Cells.Item(577, 13).Value = 10600
$ws.Cells.Item(577, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(577, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(577, 16).Value = 132
$ws.Cells.Item(577, 17).Value = 80
$ws.Cells.Item(577, 18).Value = "Hortaliza"
